$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "['MCT-2A-CAD', -]"
$ws.Range("C4").Value = "-"
$ws.Range("E4").Value = "[-, 'MCT-1A-Desenho Técnico']"
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "MCT-3A-Máquinas Térmicas e de Fluxo"
$ws.Range("D7").Value = "MCT-3A-Máquinas Térmicas e de Fluxo"
$ws.Range("E7").Value = "-"
$ws.Range("F7").Value = "-"
$ws.Range("D8").Value = "[-, 'MCT-2A-CAD']"
$ws.Range("E8").Value = "-"
